$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text, matching the
# original inline-string cell values (e.g. "252.03", "37.438.68").

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.438.68'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.052.63'
$ws.Range('E3').Value = '  +3.82%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.03'
$ws.Range('E5').Value = '  +2.69%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.654'
$ws.Range('E6').Value = '  +3.22%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '65.48'
$ws.Range('E7').Value = '  +15.15%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('E9').Value = '  +6.91%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.78'
$ws.Range('E10').Value = '  +2.34%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0768'
$ws.Range('E11').Value = '  +4.92%  '

# Row 12
$ws.Range('E12').Value = '  +1.74%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.916'
$ws.Range('E13').Value = '  -3.16%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.78'
$ws.Range('E14').Value = '  +2.66%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.352.85'
$ws.Range('E15').Value = '  +3.95%  '

# Row 16
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.26'
$ws.Range('E16').Value = '  +25.73%  '

# Row 17
$ws.Range('E17').Value = '  +6.16%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.054.71'
$ws.Range('E18').Value = '  +4.18%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '37.270.69'
$ws.Range('E19').Value = '  +5.04%  '

# Row 20
$ws.Range('E20').Value = '  +3.27%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0875'
$ws.Range('E21').Value = '  +4.28%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.45'
$ws.Range('E22').Value = '  +5.95%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.76'
$ws.Range('E23').Value = '  +3.18%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.06%  '

# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.62'
$ws.Range('E25').Value = '  +2.78%  '

# Row 26
$ws.Range('E26').Value = '  +5.08%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('E27').Value = '  +11.45%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.22'
$ws.Range('E28').Value = '  -1.16%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.99'
$ws.Range('E29').Value = '  +4.73%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.117'
$ws.Range('E30').Value = '  +24.32%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.28'
$ws.Range('E31').Value = '  +8.84%  '

# Row 32
$ws.Range('E32').Value = '  +3.59%  '

# Row 33
$ws.Range('E33').Value = '  +9.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.68'
$ws.Range('E34').Value = '  +8.92%  '

# Row 35
$ws.Range('E35').Value = '  +6.17%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  +2.59%  '

# Row 37
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.13%  '

# Row 38
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.84'
$ws.Range('E38').Value = '  +4.96%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.03'
$ws.Range('E39').Value = '  +17.14%  '

# Row 40
$ws.Range('E40').Value = '  +33.39%  '

# Row 41
$ws.Range('E41').Value = '  +18.32%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.02'
$ws.Range('E42').Value = '  +5.07%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.25'
$ws.Range('E43').Value = '  +2.76%  '

# Row 44
$ws.Range('E44').Value = '  +6.08%  '

# Row 45
$ws.Range('E45').Value = '  +4.00%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.12'
$ws.Range('E46').Value = '  +7.75%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '96.74'
$ws.Range('E47').Value = '  +6.21%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.98'
$ws.Range('E48').Value = '  +6.55%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.417.32'
$ws.Range('E49').Value = '  +3.57%  '

# Row 50
$ws.Range('E50').Value = '  +1.98%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.73'
$ws.Range('E51').Value = '  +0.93%  '
